$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Power ImpExp")

$ws.Range("D7").Value = 0.8
$ws.Range("E7").Value = 0.8

$ws.Range("E8").Select()
